$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (rows 2-101) from 4 to 5
for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 2).Value = 5
}

# Update the selected cell to C2
$ws.Range("C2").Select()
